# Fixed Non Technical skills weightage problem
# The "High Level Job Description" cell for the Business Analyst JD (E2) is
# updated to call out communication / team-player as a required soft skill,
# and the active selection is moved from F2 to E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$jdText = @'
JD For Business Analyst
Sub-Domain58; Core Banking Solution, CASA, Loan, General Ledger, Customer solution.
1.MBA Preferred (Finance &amp; Banking)
2.Should have an overall experience of minimum 3-6 years on banking domain in Product/Service Based Company and working as a Business Analyst.
3.Banking experience is Preferred.
4.Should have worked on at least one Core Banking Solution. (Ex58; Finacle, Oracle Flexcube, Fiserv, TCS BaNcs, FIS Profile etc.) 
communication , team player

Project Requirement58;
Mandatory Skill58;

1.Candidate should have experience in software development methodologies such as waterfall, agile, etc. 
2.Candidate should have experience in Gathering, analyzing &amp; documenting processes, rules and functions required to support these needs and corresponding requirements.
3.Experience in interacting with Business users &amp; SME providing recommendations to resolve issues for various business/technical groups &amp; defining strategic solutions to business problems in a multiple p
'@

$ws.Range("E2").Value = $jdText

# Move the saved selection/active cell from F2 to E2, matching the author's
# final cursor position when they committed the edit.
[void]$ws.Range("E2").Select()
